$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: text edits, in the specific order that reproduces the target
# shared-strings table layout. Editing A10 in place first reuses its
# existing shared-string slot (it becomes the new "Reiniciar..." note);
# the remaining new strings are then appended in the order they are first
# encountered while writing the other cells. ---
$ws.Range("A10").Value = "Reiniciar el variador al terminar de modificar los parametros para que los cambios tengan efecto."
$ws.Range("A2").Value = "Parametros de variador"
$ws.Range("A1").Value = "Lavadora.net"
$ws.Range("A15").Value = "El ejemplo que viene en el manual del variador esta incorrecto para el variador G5. El registro 0000h es la señal arranque/paro y 0001h es la frecuencia."
$ws.Range("A13").Value = "Mover selector de RS-422 a RS-485."
$ws.Range("A11").Value = "Reiniciar el variador al terminar de modificar los parametros para que los cambios tengan efecto."
$ws.Range("A10").ClearContents()

# --- Step 2: move the parameter table down one row (old B3:C8 -> new B4:C9) ---
$ws.Range("B3:C3").ClearContents()

$ws.Range("B4").Value = "B1-01"
$ws.Range("C4").Value = 2
$ws.Range("B5").Value = "B1-02"
$ws.Range("C5").Value = 2
$ws.Range("B6").Value = "H5-04"
$ws.Range("C6").Value = 3
$ws.Range("B7").Value = "H5-01"
$ws.Range("C7").Value = 1
$ws.Range("B8").Value = "H5-02"
$ws.Range("C8").Value = 3
$ws.Range("B9").Value = "H5-05"
$ws.Range("C9").Value = 0

# --- Step 3: sort the table alphabetically by column B ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B4:B9"))
$ws.Sort.SetRange($ws.Range("B4:C9"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# --- Step 4: border around the parameter table ---
$ws.Range("B4:C9").Borders.LineStyle = 1

# --- Step 5: "Mover selector..." note - left aligned ---
$ws.Range("A13:D13").HorizontalAlignment = -4131

# --- Step 6: the two larger notes - left/top aligned, word wrapped ---
$r11 = $ws.Range("A11:D11")
$r11.HorizontalAlignment = -4131
$r11.WrapText = $true
$r11.VerticalAlignment = -4160

$r15 = $ws.Range("A15:D15")
$r15.HorizontalAlignment = -4131
$r15.WrapText = $true
$r15.VerticalAlignment = -4160

# --- Step 7: title rows - bordered, centered ---
$r12 = $ws.Range("A1:D2")
$r12.Borders.LineStyle = 1
$r12.HorizontalAlignment = -4108

# --- Step 8: merge the four-column rows ---
$ws.Range("A1:D1").Merge()
$ws.Range("A2:D2").Merge()
$ws.Range("A11:D11").Merge()
$ws.Range("A13:D13").Merge()
$ws.Range("A15:D15").Merge()

# --- Step 9: row heights for the wrapped notes ---
$ws.Rows("11").RowHeight = 38.25
$ws.Rows("15").RowHeight = 55.5

# --- Step 10: selection / view ---
$ws.Range("A15:D15").Select()

# --- Step 11: page orientation ---
$ws.PageSetup.Orientation = 1
